$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45977.367291666669
$ws.Range("D3").Value = 45977.367291666669
$ws.Range("D4").Value = 45977.367291666669
$ws.Range("D5").Value = 45977.367291666669
$ws.Range("D6").Value = 45977.367291666669
$ws.Range("D7").Value = 45977.367291666669
$ws.Range("D8").Value = 45977.367291666669
$ws.Range("D9").Value = 45977.367291666669
$ws.Range("D10").Value = 45977.367291666669
$ws.Range("D11").Value = 45977.367291666669
$ws.Range("D12").Value = 45977.367291666669
$ws.Range("D13").Value = 45977.367291666669
$ws.Range("D14").Value = 45977.367291666669
$ws.Range("D15").Value = 45977.367291666669
$ws.Range("D16").Value = 45977.367291666669
$ws.Range("D17").Value = 45977.367291666669
$ws.Range("D18").Value = 45977.367291666669
$ws.Range("D19").Value = 45977.367291666669
$ws.Range("D20").Value = 45977.367291666669
$ws.Range("D21").Value = 45977.367291666669
$ws.Range("D22").Value = 45977.367291666669
$ws.Range("B23").Value = "103号直流"
$ws.Range("C23").Value = 45975.114421296297
$ws.Range("D23").Value = 45977.367291666669
$ws.Range("A24").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B24").Value = "101号直流"
$ws.Range("C24").Value = 45975.291817129626
$ws.Range("D24").Value = 45977.367291666669
$ws.Range("A25").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B25").Value = "905号直流"
$ws.Range("C25").Value = 45975.506874999999
$ws.Range("D25").Value = 45977.367291666669
$ws.Range("B26").Value = "502号直流"
$ws.Range("C26").Value = 45975.604062500002
$ws.Range("D26").Value = 45977.367291666669
$ws.Range("B27").Value = "201号直流"
$ws.Range("C27").Value = 45975.666678240741
$ws.Range("D27").Value = 45977.367291666669
$ws.Range("A28").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B28").Value = "604号直流"
$ws.Range("C28").Value = 45975.764236111114
$ws.Range("D28").Value = 45977.367291666669
$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "603号直流"
$ws.Range("C29").Value = 45976.044918981483
$ws.Range("D29").Value = 45977.367291666669
$ws.Range("A30").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B30").Value = "703号直流"
$ws.Range("C30").Value = 45976.053206018521
$ws.Range("D30").Value = 45977.367291666669
$ws.Range("B31").Value = "306号直流"
$ws.Range("C31").Value = 45976.109733796293
$ws.Range("D31").Value = 45977.367291666669
$ws.Range("A32").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B32").Value = "505号直流"
$ws.Range("C32").Value = 45976.327106481483
$ws.Range("D32").Value = 45977.367291666669
$ws.Range("B33").Value = "304号直流"
$ws.Range("C33").Value = 45976.376562500001
$ws.Range("D33").Value = 45977.367291666669
$ws.Range("A34").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B34").Value = "401号直流"
$ws.Range("C34").Value = 45976.434120370373
$ws.Range("D34").Value = 45977.367291666669
$ws.Range("A35").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B35").Value = "105号直流"
$ws.Range("C35").Value = 45976.513113425928
$ws.Range("D35").Value = 45977.367291666669
$ws.Range("A36").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B36").Value = "304号直流"
$ws.Range("C36").Value = 45976.538437499999
$ws.Range("D36").Value = 45977.367291666669
$ws.Range("A37").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B37").Value = "903号直流"
$ws.Range("C37").Value = 45976.542881944442
$ws.Range("D37").Value = 45977.367291666669
$ws.Range("B38").Value = "804号直流"
$ws.Range("C38").Value = 45976.544537037036
$ws.Range("D38").Value = 45977.367291666669
$ws.Range("B39").Value = "A02号直流"
$ws.Range("C39").Value = 45976.548518518517
$ws.Range("D39").Value = 45977.367291666669
$ws.Range("B40").Value = "702号直流"
$ws.Range("C40").Value = 45976.551585648151
$ws.Range("D40").Value = 45977.367291666669
$ws.Range("B41").Value = "805号直流"
$ws.Range("C41").Value = 45976.564664351848
$ws.Range("D41").Value = 45977.367291666669
$ws.Range("B42").Value = "905号直流"
$ws.Range("C42").Value = 45976.583043981482
$ws.Range("D42").Value = 45977.367291666669
$ws.Range("A43").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B43").Value = "208号直流"
$ws.Range("C43").Value = 45976.585462962961
$ws.Range("D43").Value = 45977.367291666669
$ws.Range("B44").Value = "009A号直流"
$ws.Range("C44").Value = 45976.585590277777
$ws.Range("D44").Value = 45977.367291666669
$ws.Range("A45").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B45").Value = "601号直流"
$ws.Range("C45").Value = 45976.611076388886
$ws.Range("D45").Value = 45977.367291666669
$ws.Range("A46").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B46").Value = "602号直流"
$ws.Range("C46").Value = 45976.619421296295
$ws.Range("D46").Value = 45977.367291666669
$ws.Range("A47").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B47").Value = "005B号直流"
$ws.Range("C47").Value = 45976.635196759256
$ws.Range("D47").Value = 45977.367291666669
$ws.Range("A48").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B48").Value = "007A号直流"
$ws.Range("C48").Value = 45976.654756944445
$ws.Range("D48").Value = 45977.367291666669
$ws.Range("A49").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B49").Value = "404号直流"
$ws.Range("C49").Value = 45976.67695601852
$ws.Range("D49").Value = 45977.367291666669
$ws.Range("A50").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B50").Value = "501号直流"
$ws.Range("C50").Value = 45976.690868055557
$ws.Range("D50").Value = 45977.367291666669
$ws.Range("A51").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B51").Value = "206号直流"
$ws.Range("C51").Value = 45976.699861111112
$ws.Range("D51").Value = 45977.367291666669
$ws.Range("A52").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B52").Value = "901号直流"
$ws.Range("C52").Value = 45976.780717592592
$ws.Range("D52").Value = 45977.367291666669
$ws.Range("A53").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B53").Value = "A04号直流"
$ws.Range("C53").Value = 45976.782233796293
$ws.Range("D53").Value = 45977.367291666669
$ws.Range("A54").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B54").Value = "112号直流"
$ws.Range("C54").Value = 45976.788715277777
$ws.Range("D54").Value = 45977.367291666669
$ws.Range("A55").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B55").Value = "B01号直流"
$ws.Range("C55").Value = 45976.855347222219
$ws.Range("D55").Value = 45977.367291666669
$ws.Range("A56").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B56").Value = "203号直流"
$ws.Range("C56").Value = 45976.866712962961
$ws.Range("D56").Value = 45977.367291666669

$ws.Range("E10").Select() | Out-Null
